# Inventario.xlsx edit script
# - Adds a "Precios" worksheet (purchase/sale price + profit per product)
# - Adds Precio Compra / Precio Venta / Utilidad columns to the Inventario sheet
#   (with purchase price filled in for the first four products)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Inventario sheet: new header columns D/E/F + purchase-price values
# ---------------------------------------------------------------------------
$inv = $wb.Worksheets.Item("Inventario")

# Copy the bold/centered header style from C1 onto D1:F1, then set the text
$inv.Range("C1").Copy($inv.Range("D1"))
$inv.Range("C1").Copy($inv.Range("E1"))
$inv.Range("C1").Copy($inv.Range("F1"))
$inv.Range("D1").Value = "Precio Compra"
$inv.Range("E1").Value = "Precio Venta"
$inv.Range("F1").Value = "Utilidad"

# Purchase prices for the first four products
$inv.Range("D2").Value = 8500
$inv.Range("D3").Value = 10000
$inv.Range("D4").Value = 8200
$inv.Range("D5").Value = 250

# Column widths for the new columns (target ~14.51 / ~12.56 "characters";
# the host quantizes ColumnWidth to pixel steps, so these inputs are chosen
# to land on the closest achievable stored width)
$inv.Columns.Item(4).ColumnWidth = 13.666666666666666
$inv.Columns.Item(5).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------------
# 2) New "Precios" worksheet, placed right after "Inventario"
# ---------------------------------------------------------------------------
$precios = $wb.Worksheets.Add($null, $inv)
$precios.Name = "Precios"

# Mirror column A width/style from Inventario (target ~19.66 / ~14.51 / ~12.56)
$precios.Columns.Item(1).ColumnWidth = 18.833333333333332
$precios.Columns.Item(2).ColumnWidth = 13.666666666666666
$precios.Columns.Item(3).ColumnWidth = 11.666666666666666

# Header row (bold, centered like Inventario's header)
$inv.Range("A1:D1").Copy($precios.Range("A1:D1"))
$precios.Range("A1").Value = "Producto"
$precios.Range("B1").Value = "Precio Compra"
$precios.Range("C1").Value = "Precio Venta"
$precios.Range("D1").Value = "Utilidad"

# Copy the normal product-name style from Inventario column A onto Precios column A
$inv.Range("A2:A8").Copy($precios.Range("A2:A8"))

$precios.Range("A2").Value = "Splash Victoria Secret"
$precios.Range("B2").Value = 8500
$precios.Range("C2").Value = 17000
$precios.Range("D2").Formula = "=C2-B2"

$precios.Range("A3").Value = "Splash Pink"
$precios.Range("B3").Value = 10000
$precios.Range("C3").Value = 19000
$precios.Range("D3").Formula = "=C3-B3"

$precios.Range("A4").Value = "Crema Corporal Valsy"
$precios.Range("B4").Value = 8200
$precios.Range("C4").Value = 16000
$precios.Range("D4").Formula = "=C4-B4"

$precios.Range("A5").Value = "Esencia"
$precios.Range("B5").Value = 250
$precios.Range("C5").Value = 500
$precios.Range("D5").Formula = "=C5-B5"

$precios.Range("A6").Value = "Loción Possession"
$precios.Range("B6").Value = 14000
$precios.Range("C6").Value = 32000
$precios.Range("D6").Formula = "=C6-B6"

$precios.Range("A7").Value = "Loción Expression"
$precios.Range("B7").Value = 7400
$precios.Range("C7").Value = 15000
$precios.Range("D7").Formula = "=C7-B7"

$precios.Range("A8").Value = "Loción Lotus"
$precios.Range("B8").Value = 7600
$precios.Range("C8").Value = 15000
$precios.Range("D8").Formula = "=C8-B8"

# ---------------------------------------------------------------------------
# 3) Selection / active sheet bookkeeping
# ---------------------------------------------------------------------------
[void]$inv.Range("D6").Select()
[void]$precios.Range("C12").Select()
[void]$precios.Activate()
